# ===========================================================
# Update "want-to-go" (F column) counts across sheets 1, 3, 4
# ===========================================================
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览 (Exhibitions)
$ws1.Range("F5").Value = 304
$ws1.Range("F6").Value = 393
$ws1.Range("F7").Value = 851
$ws1.Range("F9").Value = 504
$ws1.Range("F12").Value = 129
$ws1.Range("F14").Value = 229
$ws1.Range("F15").Value = 31
$ws1.Range("F16").Value = 406
$ws1.Range("F17").Value = 6581
$ws1.Range("F19").Value = 70
$ws1.Range("F21").Value = 7527
$ws1.Range("F24").Value = 3381
$ws1.Range("F26").Value = 1473
$ws1.Range("F27").Value = 881
$ws1.Range("F29").Value = 22
$ws1.Range("F30").Value = 349
$ws1.Range("F32").Value = 207
$ws1.Range("F34").Value = 1615
$ws1.Range("F36").Value = 148
$ws1.Range("F39").Value = 1187
$ws1.Range("F40").Value = 1707
$ws1.Range("F41").Value = 2129

$ws3 = $wb.Worksheets.Item(3)   # 本地生活 (Local life)
$ws3.Range("F2").Value = 240
$ws3.Range("F3").Value = 1219
$ws3.Range("F4").Value = 73

$ws4 = $wb.Worksheets.Item(4)   # 全部类型 (All types)
$ws4.Range("F2").Value = 240
$ws4.Range("F4").Value = 1219
$ws4.Range("F5").Value = 73
$ws4.Range("F7").Value = 304
$ws4.Range("F8").Value = 393
$ws4.Range("F9").Value = 851
$ws4.Range("F11").Value = 504
$ws4.Range("F15").Value = 129
$ws4.Range("F18").Value = 229
$ws4.Range("F19").Value = 31
$ws4.Range("F20").Value = 406
$ws4.Range("F21").Value = 6581
$ws4.Range("F23").Value = 70
$ws4.Range("F25").Value = 7527
$ws4.Range("F28").Value = 3381
$ws4.Range("F30").Value = 1473
$ws4.Range("F31").Value = 881
$ws4.Range("F33").Value = 22
$ws4.Range("F34").Value = 349
$ws4.Range("F37").Value = 207
$ws4.Range("F39").Value = 1615
$ws4.Range("F41").Value = 148
$ws4.Range("F44").Value = 1187
$ws4.Range("F45").Value = 1707
$ws4.Range("F47").Value = 2129
# ===========================================================
# Sheet "演出" (Performances): the 2024.04.27 "今泉爱夏" row (row 5)
# was removed. The two following events shift up one row, and the
# trailing now-duplicate row is deleted (sheet shrinks from I7 to I6).
# ===========================================================
$ws2 = $wb.Worksheets.Item(2)   # 演出 (Performances)

# --- Row 5 takes on what used to be row 6's event -----------------
# (Use NumberFormat "@" + Style reset so the literal date-looking
#  text "2024.05.18" is stored as text, not auto-converted to a date
#  serial number by the Value setter.)
$ws2.Range("B5").NumberFormat = "@"
$ws2.Range("B5").Value = "2024.05.18"
$ws2.Range("B5").Style = "Normal"
$ws2.Range("C5").Value = "北京·Rie fu日本知名唱作歌手2024出道20周年中国巡回演唱会"
$ws2.Range("D5").Value = "奥园西路1号院4-5号楼 福浪LiveHouse"
$ws2.Range("E5").Value = "2024.05.18 20:00-05.18 22:00"
$ws2.Range("F5").Value = 18
$ws2.Range("G5").Value = 380
$ws2.Range("H5").Value = "https://show.bilibili.com/platform/detail.html?id=81445"
$ws2.Range("I5").Value = "//i0.hdslb.com/bfs/openplatform/202401/6e9JD6401706239890264.jpeg"

# --- Row 6 takes on what used to be row 7's event -----------------
$ws2.Range("B6").NumberFormat = "@"
$ws2.Range("B6").Value = "2024.05.25"
$ws2.Range("B6").Style = "Normal"
$ws2.Range("C6").Value = "北京·英文原版音乐剧《剧院魅影续作：真爱永恒》Andrew Lloyd Webber’s  Love Never Dies"
$ws2.Range("D6").Value = "东直门南大街14号 北京保利剧院"
$ws2.Range("E6").Value = "2024.05.25 19:30-05.30 22:00"
$ws2.Range("F6").Value = 78
$ws2.Range("G6").Value = 680
$ws2.Range("H6").Value = "https://show.bilibili.com/platform/detail.html?id=80957"
$ws2.Range("I6").Value = "//i1.hdslb.com/bfs/openplatform/202401/0MQ6YUgo1705474811213.jpeg"

# --- The old row 7 is now a duplicate of row 6; remove it ----------
$ws2.Rows(7).Delete()
